$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.617.48"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "1.694.46"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").Value = "'312.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "'0.3946"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").Value = "'0.4028"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.06%  "
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value = "'0.9981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'53.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.33%  "
$ws.Range("D12").Value = "'0.08766"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("E13").Value = "  +12.93%  "
$ws.Range("D14").Value = "'23.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "'0.00001318"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "'7.545"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.01%  "
$ws.Range("D17").Value = "1.690.16"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "'100.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'0.07083"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.11%  "
$ws.Range("D20").Value = "'19.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").Value = "'6.690"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").Value = "24.604.01"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'3.025"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.89%  "
$ws.Range("D26").Value = "'2.312"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'22.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").Value = "'159.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "'5.160"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "'133.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "'7.548"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +31.73%  "
$ws.Range("D32").Value = "1.877.43"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("D34").Value = "'0.08644"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("D35").Value = "'7.364"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +20.85%  "
$ws.Range("D36").Value = "'1.970"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.48%  "
$ws.Range("D37").Value = "'11.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.79%  "
$ws.Range("D38").Value = "'0.2723"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").Value = "'14.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").Value = "'0.02748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.02%  "
$ws.Range("D41").Value = "'0.08980"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "'1.473"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "'0.7644"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("D44").Value = "'0.7152"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").Value = "'15.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.62%  "
$ws.Range("D46").Value = "'2.453"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").Value = "'4.161"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'139.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "'1.293"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.99%  "
$ws.Range("D51").Value = "'0.00000000379"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.21%  "
